# generate TypeTerm for Composition in domain analysis
#
# The relation "s" (Relations!row4) had source=B / target=A; this edit
# flips it to source=A / target=B, and regenerates the corresponding
# TypeTerm representation on the Terms sheet (t7: "s[B*A]" -> "s[A*B]").
# Finally, the active sheet/selection moves to the Relations sheet at
# the edited cell (D4), after having touched Terms!B10.

$wb = $excel.ActiveWorkbook

$terms = $wb.Worksheets.Item("Terms")
$relations = $wb.Worksheets.Item("Relations")

# Regenerate the TypeTerm representation for relation "s" to match its
# (now swapped) source/target concepts.
$terms.Activate()
$terms.Range("B10").Value = "s[A*B]"
$terms.Range("B10").Select()

# Swap the source/target concepts of relation "s".
$relations.Activate()
$relations.Range("C4").Value = "A"
$relations.Range("D4").Value = "B"
$relations.Range("D4").Select()
